# Re-salted the random output values (column B) for output2Salted and
# refreshed the chart's position to match the updated layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "salted" values for B1:B100 (A1:A100 stay 0..99, unchanged)
$values = @(-231,-381,461,376,-180,-264,-52,614,642,-71,190,777,1009,951,1032,1332,893,1553,1515,1695,1241,1507,1768,2426,2167,2980,3251,3240,3772,3899,3999,4310,4636,4945,4970,5373,4968,5866,5614,5943,6947,7092,7078,8062,8463,7876,9046,8709,9039,9630,10444,11156,11296,11294,11629,12014,12509,12967,13413,14063,15171,15501,15381,15700,16481,16731,17590,18467,19088,19021,20244,20994,20657,21828,22133,22496,23616,24549,24314,25659,26175,26889,26841,28220,28459,29220,29698,31092,31610,32327,32968,34004,34182,35475,35419,36684,37841,38345,38659,39844)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $values[$i]
}

# The scatter chart's anchor shifted slightly (same cell anchors, small
# offset nudge) when the workbook was regenerated - reproduce that move.
$co = $ws.ChartObjects(1)
$co.Left = 132.475
$co.Top = 9.0
